$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.283.12'
$ws.Cells.Item(2, 5).Value = '  +5.57%  '

$ws.Cells.Item(3, 4).Value = '1.916.74'
$ws.Cells.Item(3, 5).Value = '  +5.98%  '

$ws.Cells.Item(4, 4).Value = '0.9989'
$ws.Cells.Item(4, 5).Value = '  -0.11%  '

$ws.Cells.Item(5, 4).Value = '254.19'
$ws.Cells.Item(5, 5).Value = '  +1.22%  '

$ws.Cells.Item(6, 4).Value = '0.9994'
$ws.Cells.Item(6, 5).Value = '  -0.03%  '

$ws.Cells.Item(7, 4).Value = '0.5154'
$ws.Cells.Item(7, 5).Value = '  +3.76%  '

$ws.Cells.Item(8, 4).Value = '45.89'
$ws.Cells.Item(8, 5).Value = '  +7.03%  '

$ws.Cells.Item(9, 4).Value = '0.2977'
$ws.Cells.Item(9, 5).Value = '  +6.43%  '

$ws.Cells.Item(10, 4).Value = '0.06856'
$ws.Cells.Item(10, 5).Value = '  +7.44%  '

$ws.Cells.Item(11, 4).Value = '1.916.50'
$ws.Cells.Item(11, 5).Value = '  +6.02%  '

$ws.Cells.Item(12, 4).Value = '17.53'
$ws.Cells.Item(12, 5).Value = '  +4.76%  '

$ws.Cells.Item(13, 4).Value = '0.07337'
$ws.Cells.Item(13, 5).Value = '  +3.11%  '

$ws.Cells.Item(14, 4).Value = '0.6915'
$ws.Cells.Item(14, 5).Value = '  +6.75%  '

$ws.Cells.Item(15, 4).Value = '87.88'
$ws.Cells.Item(15, 5).Value = '  +7.48%  '

$ws.Cells.Item(16, 4).Value = '4.921'
$ws.Cells.Item(16, 5).Value = '  +4.40%  '

$ws.Cells.Item(17, 4).Value = '30.278.90'
$ws.Cells.Item(17, 5).Value = '  +5.60%  '

$ws.Cells.Item(18, 4).Value = '0.000008151'
$ws.Cells.Item(18, 5).Value = '  +10.73%  '

$ws.Cells.Item(19, 4).Value = '0.9998'
$ws.Cells.Item(19, 5).Value = '  +0.05%  '

$ws.Cells.Item(20, 4).Value = '13.07'
$ws.Cells.Item(20, 5).Value = '  +6.59%  '

$ws.Cells.Item(21, 4).Value = '2.164.46'
$ws.Cells.Item(21, 5).Value = '  +6.27%  '

$ws.Cells.Item(22, 4).Value = '0.9981'
$ws.Cells.Item(22, 5).Value = '  -0.21%  '

$ws.Cells.Item(23, 4).Value = '4.865'
$ws.Cells.Item(23, 5).Value = '  +5.14%  '

$ws.Cells.Item(24, 4).Value = '5.760'
$ws.Cells.Item(24, 5).Value = '  +8.46%  '

$ws.Cells.Item(25, 4).Value = '9.195'
$ws.Cells.Item(25, 5).Value = '  +3.34%  '

$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(26, 4).Value = '146.73'
$ws.Cells.Item(26, 5).Value = '  +2.92%  '

$ws.Cells.Item(27, 2).Value = 'BitcoinCash'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(27, 4).Value = '139.47'
$ws.Cells.Item(27, 5).Value = '  +23.50%  '

$ws.Cells.Item(28, 5).Value = '  +7.81%  '

$ws.Cells.Item(29, 4).Value = '2.011'
$ws.Cells.Item(29, 5).Value = '  +6.73%  '

$ws.Cells.Item(30, 4).Value = '1.380'
$ws.Cells.Item(30, 5).Value = '  -1.49%  '

$ws.Cells.Item(31, 4).Value = '4.287'
$ws.Cells.Item(31, 5).Value = '  +2.62%  '

$ws.Cells.Item(32, 4).Value = '0.08860'
$ws.Cells.Item(32, 5).Value = '  +5.90%  '

$ws.Cells.Item(33, 4).Value = '4.033'
$ws.Cells.Item(33, 5).Value = '  +5.18%  '

$ws.Cells.Item(34, 4).Value = '0.05141'
$ws.Cells.Item(34, 5).Value = '  +3.59%  '

$ws.Cells.Item(35, 5).Value = '  +6.71%  '

$ws.Cells.Item(36, 4).Value = '0.7209'
$ws.Cells.Item(36, 5).Value = '  +6.80%  '

$ws.Cells.Item(37, 4).Value = '2.684'
$ws.Cells.Item(37, 5).Value = '  +0.67%  '

$ws.Cells.Item(38, 4).Value = '2.308'
$ws.Cells.Item(38, 5).Value = '  +7.74%  '

$ws.Cells.Item(39, 4).Value = '2.825'
$ws.Cells.Item(39, 5).Value = '  +5.82%  '

$ws.Cells.Item(40, 4).Value = '0.9782'
$ws.Cells.Item(40, 5).Value = '  +1.85%  '

$ws.Cells.Item(41, 4).Value = '0.01704'
$ws.Cells.Item(41, 5).Value = '  +7.22%  '

$ws.Cells.Item(42, 4).Value = '6.128'
$ws.Cells.Item(42, 5).Value = '  +3.34%  '

$ws.Cells.Item(43, 2).Value = 'Quant'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(43, 4).Value = '106.27'
$ws.Cells.Item(43, 5).Value = '  +5.23%  '

$ws.Cells.Item(44, 2).Value = 'TheSandbox'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(44, 4).Value = '0.4327'
$ws.Cells.Item(44, 5).Value = '  +5.11%  '

$ws.Cells.Item(45, 4).Value = '0.9990'
$ws.Cells.Item(45, 5).Value = '  -0.10%  '

$ws.Cells.Item(46, 4).Value = '7.722'
$ws.Cells.Item(46, 5).Value = '  +7.13%  '

$ws.Cells.Item(47, 4).Value = '0.1279'
$ws.Cells.Item(47, 5).Value = '  +4.58%  '

$ws.Cells.Item(48, 4).Value = '0.05738'
$ws.Cells.Item(48, 5).Value = '  +4.42%  '

$ws.Cells.Item(49, 2).Value = 'Elrond'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(49, 4).Value = '33.45'
$ws.Cells.Item(49, 5).Value = '  +6.51%  '

$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = '8.535'
$ws.Cells.Item(50, 5).Value = '  +4.09%  '

$ws.Cells.Item(51, 4).Value = '0.3845'
$ws.Cells.Item(51, 5).Value = '  +6.34%  '
